$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "28.017.20"
$ws.Range("E2").Value = "  -1.06%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.762.44"
$ws.Range("E3").Value = "  -1.54%  "
$ws.Range("E4").Value = "  +0.06%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "335.55"
$ws.Range("E5").Value = "  -1.16%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.9996"
$ws.Range("E6").Value = "  +0.07%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.3926"
$ws.Range("E7").Value = "  +2.36%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3394"
$ws.Range("E8").Value = "  -1.72%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "45.38"
$ws.Range("E9").Value = "  -3.60%  "
$ws.Range("E10").Value = "  -3.08%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07239"
$ws.Range("E11").Value = "  -2.52%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.001"
$ws.Range("E12").Value = "  +0.14%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "22.35"
$ws.Range("E13").Value = "  -4.16%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "6.154"
$ws.Range("E14").Value = "  -5.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.114"
$ws.Range("E15").Value = "  -3.40%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "1.759.38"
$ws.Range("E16").Value = "  -1.43%  "
$ws.Range("E17").Value = "  -2.05%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.06631"
$ws.Range("E18").Value = "  -0.80%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "80.47"
$ws.Range("E19").Value = "  -2.57%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.9994"
$ws.Range("E20").Value = "  -0.01%  "
$ws.Range("E21").Value = "  -3.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.231"
$ws.Range("E22").Value = "  -3.93%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "28.010.34"
$ws.Range("E23").Value = "  -1.06%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "11.66"
$ws.Range("E24").Value = "  -3.83%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.390"
$ws.Range("E25").Value = "  +1.60%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "155.00"
$ws.Range("E26").Value = "  +0.03%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "19.96"
$ws.Range("E27").Value = "  -4.08%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "2.318"
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.960.05"
$ws.Range("E29").Value = "  -1.42%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.280"
$ws.Range("E30").Value = "  -11.71%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "129.43"
$ws.Range("E31").Value = "  -4.85%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.081"
$ws.Range("E32").Value = "  +2.80%  "
$ws.Range("E33").Value = "  -5.22%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.08726"
$ws.Range("E34").Value = "  -2.45%  "
$ws.Range("E35").Value = "  -5.85%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.06182"
$ws.Range("E36").Value = "  -3.50%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.02291"
$ws.Range("E37").Value = "  -5.86%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.147"
$ws.Range("E38").Value = "  -4.13%  "
$ws.Range("E39").Value = "  -5.55%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.2115"
$ws.Range("E40").Value = "  -2.95%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.499"
$ws.Range("E41").Value = "  -0.24%  "
$ws.Range("E42").Value = "  -3.52%  "
$ws.Range("E43").Value = "  +0.07%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "7.891"
$ws.Range("E44").Value = "  -5.14%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.72"
$ws.Range("E45").Value = "  -3.97%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "3.833"
$ws.Range("E46").Value = "  -1.20%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.6000"
$ws.Range("E47").Value = "  -5.24%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "126.96"
$ws.Range("E48").Value = "  -5.07%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.990"
$ws.Range("E49").Value = "  -4.86%  "
$ws.Range("E50").Value = "  -4.62%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.07003"
$ws.Range("E51").Value = "  -6.57%  "
